$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-18 18:33:11"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
